$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 250003000
$ws.Range("J51").Value = 6000
$ws.Range("L51").Value = 6000
$ws.Range("N51").Value = -6968
$ws.Range("H53").Value = 607.8889
$ws.Range("I53").Value = 613.9167
$ws.Range("K53").Value = 613.9167
$ws.Range("M53").Value = 23.08330000000001
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("H64").Value = 125005304
$ws.Range("J64").Value = 333337820
$ws.Range("L64").Value = 333337820
$ws.Range("N64").Value = -333338316
$ws.Range("H67").Value = 125005304
$ws.Range("J67").Value = 333337820
$ws.Range("L67").Value = 333337820
$ws.Range("N67").Value = -333339536
$ws.Range("H88").Value = 3580274.5
$ws.Range("I88").Value = 9102337
$ws.Range("J88").Value = 7175
$ws.Range("K88").Value = 9102337
$ws.Range("L88").Value = 7175
$ws.Range("M88").Value = -9101931
$ws.Range("N88").Value = -7987
$ws.Range("H91").Value = 3580274.5
$ws.Range("I91").Value = 9102337
$ws.Range("J91").Value = 7175
$ws.Range("K91").Value = 9102337
$ws.Range("L91").Value = 7175
$ws.Range("M91").Value = -9100933
$ws.Range("N91").Value = -9983
$ws.Range("H98").Value = 2264.1482
$ws.Range("I98").Value = 2264.1482
$ws.Range("K98").Value = 2264.1482
$ws.Range("M98").Value = -766.1482000000001
$ws.Range("H107").Value = 410.33334
$ws.Range("J107").Value = 392.2
$ws.Range("L107").Value = 392.2
$ws.Range("N107").Value = -4232.2
$ws.Range("H111").Value = 1374.3334
$ws.Range("I111").Value = 894.8333
$ws.Range("J111").Value = 2333.3333
$ws.Range("K111").Value = 2684.4999
$ws.Range("L111").Value = 6999.999899999999
$ws.Range("M111").Value = 382.5001000000002
$ws.Range("N111").Value = -13133.9999
$ws.Range("H118").Value = 743.5
$ws.Range("I118").Value = 743.5
$ws.Range("K118").Value = 2230.5
$ws.Range("M118").Value = -573.5
$ws.Range("H122").Value = 2264.1482
$ws.Range("I122").Value = 2264.1482
$ws.Range("K122").Value = 6792.444600000001
$ws.Range("M122").Value = -4342.444600000001
$ws.Range("H129").Value = 2422.7334
$ws.Range("I129").Value = 667.75
$ws.Range("J129").Value = 4428.4287
$ws.Range("K129").Value = 2003.25
$ws.Range("L129").Value = 13285.2861
$ws.Range("M129").Value = 2996.75
$ws.Range("N129").Value = -23285.2861
$ws.Range("H137").Value = 1516831.2
$ws.Range("I137").Value = 2001509.1
$ws.Range("J137").Value = 2213
$ws.Range("K137").Value = 6004527.300000001
$ws.Range("L137").Value = 6639
$ws.Range("M137").Value = -6001977.300000001
$ws.Range("N137").Value = -11739
$ws.Range("H138").Value = 4254.256
$ws.Range("J138").Value = 3846.9443
$ws.Range("L138").Value = 11540.8329
$ws.Range("N138").Value = -21820.8329
$ws.Range("M61").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 1000
$ws.Range("M8").Value = -856
$ws.Range("H13").Value = 299
$ws.Range("I13").Value = 299
$ws.Range("K13").Value = 299
$ws.Range("M13").Value = -155
$ws.Range("H14").Value = 7420.375
$ws.Range("J14").Value = 364.25
$ws.Range("L14").Value = 364.25
$ws.Range("N14").Value = -714.25
$ws.Range("H19").Value = 504
$ws.Range("I19").Value = 504
$ws.Range("K19").Value = 504
$ws.Range("M19").Value = -275
$ws.Range("H21").Value = 1572
$ws.Range("J21").Value = 1150
$ws.Range("L21").Value = 1150
$ws.Range("N21").Value = -1898
$ws.Range("H29").Value = 950
$ws.Range("J29").Value = 1400
$ws.Range("L29").Value = 1400
$ws.Range("N29").Value = -2016
$ws.Range("H30").Value = 738.625
$ws.Range("J30").Value = 1000
$ws.Range("L30").Value = 1000
$ws.Range("N30").Value = -1300
$ws.Range("H32").Value = 2250.3784
$ws.Range("I32").Value = 2327.6287
$ws.Range("J32").Value = 898.5
$ws.Range("K32").Value = 2327.6287
$ws.Range("L32").Value = 898.5
$ws.Range("M32").Value = -2040.6287
$ws.Range("N32").Value = -1472.5
$ws.Range("H34").Value = 35000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 35000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 35000
$ws.Range("N34").Value = -35542
$ws.Range("H36").Value = 3271.2856
$ws.Range("I36").Value = 4966.3335
$ws.Range("K36").Value = 4966.3335
$ws.Range("M36").Value = -4620.3335
$ws.Range("H37").Value = 40845
$ws.Range("I37").Value = 36677.668
$ws.Range("J37").Value = 45012.332
$ws.Range("K37").Value = 36677.668
$ws.Range("L37").Value = 45012.332
$ws.Range("M37").Value = -36404.668
$ws.Range("N37").Value = -45558.332
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("H45").Value = 30043
$ws.Range("J45").Value = 2205.6667
$ws.Range("L45").Value = 2205.6667
$ws.Range("N45").Value = -2959.6667
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 30000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31538
$ws.Range("H55").Value = 77053
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("H61").Value = 3684.8572
$ws.Range("J61").Value = 3882
$ws.Range("L61").Value = 3882
$ws.Range("N61").Value = -4306
$ws.Range("H74").Value = 187064.03
$ws.Range("I74").Value = 254196.95
$ws.Range("K74").Value = 254196.95
$ws.Range("M74").Value = -253322.95
$ws.Range("H77").Value = 187064.03
$ws.Range("I77").Value = 254196.95
$ws.Range("K77").Value = 1270984.75
$ws.Range("M77").Value = -1266616.75
$ws.Range("H97").Value = 928.2273
$ws.Range("I97").Value = 917.9722
$ws.Range("K97").Value = 917.9722
$ws.Range("M97").Value = -421.9722
$ws.Range("H110").Value = 833.3333
$ws.Range("I110").Value = 850
$ws.Range("K110").Value = 850
$ws.Range("M110").Value = 1195
$ws.Range("H122").Value = 5717.448
$ws.Range("I122").Value = 5851.72
$ws.Range("J122").Value = 4878.25
$ws.Range("K122").Value = 17555.16
$ws.Range("L122").Value = 14634.75
$ws.Range("M122").Value = -15105.16
$ws.Range("N122").Value = -19534.75
$ws.Range("H136").Value = 3684.8572
$ws.Range("J136").Value = 3882
$ws.Range("L136").Value = 11646
$ws.Range("N136").Value = -16746
$ws.Range("M34").ClearContents()
$ws.Range("M42").ClearContents()
$ws.Range("M54").ClearContents()
$ws.Range("M55").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13002514
$ws.Range("I105").Value = 1113944
$ws.Range("J105").Value = 22729526
$ws.Range("K105").Value = 1113944
$ws.Range("L105").Value = 22729526
$ws.Range("M105").Value = -1112197
$ws.Range("N105").Value = -22733020
$ws.Range("H107").Value = 2138060.5
$ws.Range("I107").Value = 2850250.2
$ws.Range("J107").Value = 1491.2222
$ws.Range("K107").Value = 2850250.2
$ws.Range("L107").Value = 1491.2222
$ws.Range("M107").Value = -2848330.2
$ws.Range("N107").Value = -5331.2222
$ws.Range("H134").Value = 2819.1738
$ws.Range("I134").Value = 2662.7334
$ws.Range("K134").Value = 7988.2002
$ws.Range("M134").Value = -5453.2002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 265.14285
$ws.Range("I10").Value = 276
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 276
$ws.Range("L10").Value = 200
$ws.Range("M10").Value = -137
$ws.Range("N10").Value = -478
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("H16").Value = 2045
$ws.Range("I16").Value = 2049
$ws.Range("K16").Value = 2049
$ws.Range("M16").Value = -1762
$ws.Range("H21").Value = 3999.5
$ws.Range("I21").Value = 3999
$ws.Range("K21").Value = 3999
$ws.Range("M21").Value = -3764
$ws.Range("H22").Value = 998.05554
$ws.Range("I22").Value = 1044
$ws.Range("J22").Value = 906.1667
$ws.Range("K22").Value = 1044
$ws.Range("L22").Value = 906.1667
$ws.Range("M22").Value = -694
$ws.Range("N22").Value = -1606.1667
$ws.Range("H23").Value = 5002499.5
$ws.Range("I23").Value = 5002499.5
$ws.Range("K23").Value = 5002499.5
$ws.Range("M23").Value = -5002259.5
$ws.Range("H27").Value = 5002499.5
$ws.Range("I27").Value = 5002499.5
$ws.Range("K27").Value = 5002499.5
$ws.Range("M27").Value = -5002307.5
$ws.Range("H31").Value = 4634994
$ws.Range("I31").Value = 5113.875
$ws.Range("J31").Value = 11369365
$ws.Range("K31").Value = 5113.875
$ws.Range("L31").Value = 11369365
$ws.Range("M31").Value = -4818.875
$ws.Range("N31").Value = -11369955
$ws.Range("H34").Value = 4634994
$ws.Range("I34").Value = 5113.875
$ws.Range("J34").Value = 11369365
$ws.Range("K34").Value = 5113.875
$ws.Range("L34").Value = 11369365
$ws.Range("M34").Value = -4911.875
$ws.Range("N34").Value = -11369769
$ws.Range("H58").Value = 1409.5151
$ws.Range("J58").Value = 2799.7144
$ws.Range("L58").Value = 2799.7144
$ws.Range("N58").Value = -3205.7144
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2
$ws.Range("H107").Value = 3572291.8
$ws.Range("I107").Value = 7143284.5
$ws.Range("J107").Value = 1299
$ws.Range("K107").Value = 7143284.5
$ws.Range("L107").Value = 1299
$ws.Range("M107").Value = -7141364.5
$ws.Range("N107").Value = -5139
$ws.Range("H113").Value = 2045
$ws.Range("I113").Value = 2049
$ws.Range("K113").Value = 2049
$ws.Range("M113").Value = 121
$ws.Range("H122").Value = 4347.1875
$ws.Range("I122").Value = 3509.375
$ws.Range("K122").Value = 10528.125
$ws.Range("M122").Value = -8078.125
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H132").Value = 11113113
$ws.Range("J132").Value = 41669880
$ws.Range("L132").Value = 125009640
$ws.Range("N132").Value = -125014700
$ws.Range("H134").Value = 2223.6
$ws.Range("I134").Value = 1951
$ws.Range("J134").Value = 3314
$ws.Range("K134").Value = 5853
$ws.Range("L134").Value = 9942
$ws.Range("M134").Value = -3318
$ws.Range("N134").Value = -15012
$ws.Range("H136").Value = 1409.5151
$ws.Range("J136").Value = 2799.7144
$ws.Range("L136").Value = 8399.143199999999
$ws.Range("N136").Value = -13499.1432
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1001601.6
$ws.Range("J97").Value = 2002
$ws.Range("L97").Value = 6006
$ws.Range("N97").Value = -6998
$ws.Range("H134").Value = 1738.3529
$ws.Range("I134").Value = 1176
$ws.Range("K134").Value = 3528
$ws.Range("M134").Value = 1542

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 8333.333000000001
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("H8").Value = 8333.333000000001
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("H19").Value = 829.6
$ws.Range("I19").Value = 983
$ws.Range("K19").Value = 983
$ws.Range("M19").Value = -695
$ws.Range("H24").Value = 10118.25
$ws.Range("I24").Value = 10996.5
$ws.Range("J24").Value = 9240
$ws.Range("K24").Value = 10996.5
$ws.Range("L24").Value = 9240
$ws.Range("M24").Value = -10823.5
$ws.Range("N24").Value = -9586
$ws.Range("H70").Value = 92451.39
$ws.Range("I70").Value = 130180.5
$ws.Range("J70").Value = 6213.4287
$ws.Range("K70").Value = 130180.5
$ws.Range("L70").Value = 6213.4287
$ws.Range("M70").Value = -129910.5
$ws.Range("N70").Value = -6753.4287
$ws.Range("H73").Value = 92451.39
$ws.Range("I73").Value = 130180.5
$ws.Range("J73").Value = 6213.4287
$ws.Range("K73").Value = 130180.5
$ws.Range("L73").Value = 6213.4287
$ws.Range("M73").Value = -129244.5
$ws.Range("N73").Value = -8085.4287
$ws.Range("H113").Value = 3236.8462
$ws.Range("I113").Value = 3772.7144
$ws.Range("J113").Value = 2611.6667
$ws.Range("K113").Value = 3772.7144
$ws.Range("L113").Value = 2611.6667
$ws.Range("M113").Value = -1602.7144
$ws.Range("N113").Value = -6951.6667
$ws.Range("H122").Value = 4196.1665
$ws.Range("I122").Value = 3668.5454
$ws.Range("K122").Value = 11005.6362
$ws.Range("M122").Value = -8555.636200000001
$ws.Range("H126").Value = 15871.5
$ws.Range("I126").Value = 2633.7778
$ws.Range("J126").Value = 39699.4
$ws.Range("K126").Value = 7901.3334
$ws.Range("L126").Value = 119098.2
$ws.Range("M126").Value = -5431.3334
$ws.Range("N126").Value = -124038.2
$ws.Range("H132").Value = 2999.7144
$ws.Range("I132").Value = 2599.7
$ws.Range("J132").Value = 3999.75
$ws.Range("K132").Value = 7799.099999999999
$ws.Range("L132").Value = 11999.25
$ws.Range("M132").Value = -5269.099999999999
$ws.Range("N132").Value = -17059.25
$ws.Range("M7").ClearContents()
$ws.Range("M8").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6005.8335
$ws.Range("I7").Value = 4894.5557
$ws.Range("J7").Value = 7117.1113
$ws.Range("K7").Value = 4894.5557
$ws.Range("L7").Value = 7117.1113
$ws.Range("M7").Value = -4782.5557
$ws.Range("N7").Value = -7341.1113
$ws.Range("H12").Value = 3763.8
$ws.Range("J12").Value = 4580
$ws.Range("L12").Value = 4580
$ws.Range("N12").Value = -4920
$ws.Range("H17").Value = 436.5
$ws.Range("I17").Value = 436.5
$ws.Range("K17").Value = 436.5
$ws.Range("M17").Value = -266.5
$ws.Range("H24").Value = 4999.5
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("H25").Value = 17999
$ws.Range("I25").Value = 17999
$ws.Range("K25").Value = 17999
$ws.Range("M25").Value = -17769
$ws.Range("H40").Value = 4436.778
$ws.Range("I40").Value = 1822
$ws.Range("J40").Value = 9666.333000000001
$ws.Range("K40").Value = 1822
$ws.Range("L40").Value = 9666.333000000001
$ws.Range("M40").Value = -1686
$ws.Range("N40").Value = -9938.333000000001
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("H56").Value = 14002.667
$ws.Range("I56").Value = 14002.667
$ws.Range("K56").Value = 14002.667
$ws.Range("M56").Value = -13311.667
$ws.Range("H122").Value = 6553.88
$ws.Range("I122").Value = 5265.864
$ws.Range("J122").Value = 15999.333
$ws.Range("K122").Value = 15797.592
$ws.Range("L122").Value = 47997.999
$ws.Range("M122").Value = -13347.592
$ws.Range("N122").Value = -52897.999
$ws.Range("H126").Value = 6005.8335
$ws.Range("I126").Value = 4894.5557
$ws.Range("J126").Value = 7117.1113
$ws.Range("K126").Value = 14683.6671
$ws.Range("L126").Value = 21351.3339
$ws.Range("M126").Value = -12213.6671
$ws.Range("N126").Value = -26291.3339
$ws.Range("H132").Value = 9314.929
$ws.Range("I132").Value = 6157.778
$ws.Range("J132").Value = 14997.8
$ws.Range("K132").Value = 18473.334
$ws.Range("L132").Value = 44993.39999999999
$ws.Range("M132").Value = -15943.334
$ws.Range("N132").Value = -50053.39999999999
$ws.Range("H136").Value = 4909.8623
$ws.Range("I136").Value = 3853.625
$ws.Range("J136").Value = 9979.799999999999
$ws.Range("K136").Value = 11560.875
$ws.Range("L136").Value = 29939.4
$ws.Range("M136").Value = -9010.875
$ws.Range("N136").Value = -35039.39999999999
$ws.Range("M24").ClearContents()
$ws.Range("M48").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 11495
$ws.Range("I3").Value = 7990
$ws.Range("K3").Value = 7990
$ws.Range("M3").Value = -7876
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("H10").Value = 24999.5
$ws.Range("J10").Value = 25000
$ws.Range("L10").Value = 25000
$ws.Range("N10").Value = -25338
$ws.Range("H12").Value = 9666.333000000001
$ws.Range("I12").Value = 4999
$ws.Range("K12").Value = 4999
$ws.Range("M12").Value = -4857
$ws.Range("H13").Value = 5311.6665
$ws.Range("J13").Value = 5311.6665
$ws.Range("L13").Value = 5311.6665
$ws.Range("N13").Value = -5591.6665
$ws.Range("H17").Value = 5682.25
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H23").Value = 26281.8
$ws.Range("J23").Value = 14949.5
$ws.Range("L23").Value = 14949.5
$ws.Range("N23").Value = -15407.5
$ws.Range("H30").Value = 30009
$ws.Range("I30").Value = 30009
$ws.Range("K30").Value = 30009
$ws.Range("M30").Value = -29902
$ws.Range("H31").Value = 8997.5
$ws.Range("I31").Value = 10000
$ws.Range("J31").Value = 7995
$ws.Range("K31").Value = 10000
$ws.Range("L31").Value = 7995
$ws.Range("M31").Value = -9652
$ws.Range("N31").Value = -8691
$ws.Range("H34").Value = 15000
$ws.Range("I34").Value = 15000
$ws.Range("K34").Value = 15000
$ws.Range("M34").Value = -14797
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("H42").Value = 710000
$ws.Range("I42").Value = 710000
$ws.Range("K42").Value = 710000
$ws.Range("M42").Value = -709622
$ws.Range("H43").Value = 42465
$ws.Range("I43").Value = 30000
$ws.Range("K43").Value = 30000
$ws.Range("M43").Value = -29851
$ws.Range("H45").Value = 10544.667
$ws.Range("I45").Value = 3000
$ws.Range("J45").Value = 11230.546
$ws.Range("K45").Value = 3000
$ws.Range("L45").Value = 11230.546
$ws.Range("M45").Value = -2509
$ws.Range("N45").Value = -12212.546
$ws.Range("H54").Value = 49988.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 49988.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 49988.5
$ws.Range("N54").Value = -51028.5
$ws.Range("H58").Value = 11493.5
$ws.Range("I58").Value = 11493.5
$ws.Range("K58").Value = 11493.5
$ws.Range("M58").Value = -11185.5
$ws.Range("H81").Value = 4470.2085
$ws.Range("I81").Value = 4249.1875
$ws.Range("J81").Value = 4912.25
$ws.Range("K81").Value = 8498.375
$ws.Range("L81").Value = 9824.5
$ws.Range("M81").Value = -7437.375
$ws.Range("N81").Value = -11946.5
$ws.Range("H84").Value = 4470.2085
$ws.Range("I84").Value = 4249.1875
$ws.Range("J84").Value = 4912.25
$ws.Range("K84").Value = 42491.875
$ws.Range("L84").Value = 49122.5
$ws.Range("M84").Value = -37187.875
$ws.Range("N84").Value = -59730.5
$ws.Range("H100").Value = 38462076
$ws.Range("I100").Value = 441.3889
$ws.Range("K100").Value = 882.7778
$ws.Range("M100").Value = -341.7778
$ws.Range("H107").Value = 1236.091
$ws.Range("J107").Value = 1113.8572
$ws.Range("L107").Value = 3341.5716
$ws.Range("N107").Value = -7181.571599999999
$ws.Range("H122").Value = 16668870
$ws.Range("I122").Value = 2385.4167
$ws.Range("K122").Value = 7156.250100000001
$ws.Range("M122").Value = -4706.250100000001
$ws.Range("H126").Value = 6497.4585
$ws.Range("I126").Value = 7582.4736
$ws.Range("J126").Value = 2374.4
$ws.Range("K126").Value = 22747.4208
$ws.Range("L126").Value = 7123.200000000001
$ws.Range("M126").Value = -20277.4208
$ws.Range("N126").Value = -12063.2
$ws.Range("H132").Value = 5278.3335
$ws.Range("I132").Value = 4994.1
$ws.Range("K132").Value = 14982.3
$ws.Range("M132").Value = -12452.3
$ws.Range("M7").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("M35").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("M54").ClearContents()
